$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost (Q2) and Nord (R2) coordinate values to whole numbers
$ws.Range("Q2").Value2 = [math]::Round([double]$ws.Range("Q2").Value2, 0)
$ws.Range("R2").Value2 = [math]::Round([double]$ws.Range("R2").Value2, 0)

# Clear the Starttid (Z2) and Sluttid (AB2) cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
